$wb = $excel.ActiveWorkbook

# ---- Sheet "2018": add RATING / Slow row (row 8) ----
$ws1 = $wb.Worksheets.Item("2018")

$ws1.Range("A8").Value = "RATING"
$f1 = $ws1.Range("A8").Font
$f1.Name = "Arial"
$f1.Size = 12
$f1.Bold = $true

$ws1.Range("C8").Value = "Slow"
$ws1.Range("C8").HorizontalAlignment = -4108
$f2 = $ws1.Range("C8").Font
$f2.Name = "Arial"
$f2.Size = 12
$f2.Bold = $true
$f2.Color = 15773696

$ws1.Rows.Item(8).RowHeight = 15.75
$ws1.Range("C8").Select()

# ---- Sheet "Overall Stats": add RATING / Slow row (row 11) ----
$ws2 = $wb.Worksheets.Item("Overall Stats")

$ws2.Range("A11").Value = "RATING"
$f3 = $ws2.Range("A11").Font
$f3.Name = "Arial"
$f3.Size = 12
$f3.Bold = $true

$ws2.Range("B11").Value = "Slow"
$ws2.Range("B11").HorizontalAlignment = -4152
$f4 = $ws2.Range("B11").Font
$f4.Name = "Arial"
$f4.Size = 12
$f4.Bold = $true
$f4.Color = 15773696

$ws2.Rows.Item(11).RowHeight = 15.75
$ws2.Range("A27").Select()
